$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source price column stores values as literal text (e.g. '18.24',
# '20.472.41'); pre-format the specific cells whose new price would otherwise
# be auto-parsed as a number so they keep their original text representation.
# (NumberFormat is applied per-cell because setting it on a multi-area union
# range only affects the first area.)
$textPriceCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D16", "D18", "D19", "D20", "D21", "D22", "D24", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "20.478.85"
$ws.Range("E2").Value = "  +2.75%  "

# Row 3
$ws.Range("D3").Value = "1.472.02"
$ws.Range("E3").Value = "  +4.16%  "

# Row 4
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +0.54%  "

# Row 5
$ws.Range("D5").Value = "0.9665"
$ws.Range("E5").Value = "  -3.55%  "

# Row 6
$ws.Range("D6").Value = "275.74"
$ws.Range("E6").Value = "  -0.20%  "

# Row 7
$ws.Range("D7").Value = "0.3651"
$ws.Range("E7").Value = "  -1.10%  "

# Row 8
$ws.Range("D8").Value = "0.3065"
$ws.Range("E8").Value = "  -1.26%  "

# Row 9
$ws.Range("D9").Value = "39.86"
$ws.Range("E9").Value = "  -0.07%  "

# Row 10
$ws.Range("D10").Value = "1.047"
$ws.Range("E10").Value = "  +1.28%  "

# Row 11
$ws.Range("D11").Value = "0.06613"
$ws.Range("E11").Value = "  +1.57%  "

# Row 12
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").Value = "  -0.04%  "

# Row 13
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").Value = "18.24"
$ws.Range("E13").Value = "  +3.76%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.455"
$ws.Range("E14").Value = "  -0.18%  "

# Row 15
$ws.Range("E15").Value = "  -0.48%  "

# Row 16
$ws.Range("D16").Value = "0.00001029"
$ws.Range("E16").Value = "  +1.02%  "

# Row 17
$ws.Range("D17").Value = "1.472.41"
$ws.Range("E17").Value = "  +3.92%  "

# Row 18
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "0.05898"
$ws.Range("E18").Value = "  +3.68%  "

# Row 19
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "0.9739"
$ws.Range("E19").Value = "  -2.81%  "

# Row 20
$ws.Range("D20").Value = "69.07"
$ws.Range("E20").Value = "  -2.55%  "

# Row 21
$ws.Range("D21").Value = "5.448"
$ws.Range("E21").Value = "  -2.58%  "

# Row 22
$ws.Range("D22").Value = "14.40"
$ws.Range("E22").Value = "  -2.18%  "

# Row 23
$ws.Range("E23").Value = "  -0.12%  "

# Row 24
$ws.Range("D24").Value = "2.249"
$ws.Range("E24").Value = "  +0.74%  "

# Row 25
$ws.Range("D25").Value = "20.524.08"
$ws.Range("E25").Value = "  +2.79%  "

# Row 26
$ws.Range("D26").Value = "141.77"
$ws.Range("E26").Value = "  +6.58%  "

# Row 27
$ws.Range("D27").Value = "2.135"
$ws.Range("E27").Value = "  -5.72%  "

# Row 28
$ws.Range("D28").Value = "17.22"
$ws.Range("E28").Value = "  +0.27%  "

# Row 29
$ws.Range("D29").Value = "1.627.83"
$ws.Range("E29").Value = "  +3.28%  "

# Row 30
$ws.Range("D30").Value = "113.67"
$ws.Range("E30").Value = "  +3.35%  "

# Row 31
$ws.Range("D31").Value = "3.880"
$ws.Range("E31").Value = "  -1.02%  "

# Row 32
$ws.Range("D32").Value = "4.953"
$ws.Range("E32").Value = "  -4.68%  "

# Row 33
$ws.Range("D33").Value = "0.8031"
$ws.Range("E33").Value = "  -0.62%  "

# Row 34
$ws.Range("D34").Value = "0.07873"
$ws.Range("E34").Value = "  +1.36%  "

# Row 35
$ws.Range("D35").Value = "1.529"
$ws.Range("E35").Value = "  +3.47%  "

# Row 36
$ws.Range("D36").Value = "1.239"
$ws.Range("E36").Value = "  +12.18%  "

# Row 37
$ws.Range("D37").Value = "0.05730"
$ws.Range("E37").Value = "  -1.72%  "

# Row 38
$ws.Range("D38").Value = "4.752"
$ws.Range("E38").Value = "  -2.64%  "

# Row 39
$ws.Range("D39").Value = "0.9697"
$ws.Range("E39").Value = "  -3.13%  "

# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.02036"
$ws.Range("E40").Value = "  -0.29%  "

# Row 41
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "7.643"
$ws.Range("E41").Value = "  -6.01%  "

# Row 42
$ws.Range("D42").Value = "10.42"
$ws.Range("E42").Value = "  -0.02%  "

# Row 43
$ws.Range("D43").Value = "0.1879"
$ws.Range("E43").Value = "  -0.24%  "

# Row 44
$ws.Range("D44").Value = "0.5284"
$ws.Range("E44").Value = "  -0.27%  "

# Row 45
$ws.Range("D45").Value = "3.502"
$ws.Range("E45").Value = "  -0.89%  "

# Row 46
$ws.Range("D46").Value = "12.05"
$ws.Range("E46").Value = "  -2.30%  "

# Row 47
$ws.Range("D47").Value = "117.01"
$ws.Range("E47").Value = "  +0.44%  "

# Row 48
$ws.Range("D48").Value = "0.5169"
$ws.Range("E48").Value = "  -0.12%  "

# Row 49
$ws.Range("D49").Value = "1.768"
$ws.Range("E49").Value = "  +0.19%  "

# Row 50
$ws.Range("D50").Value = "0.06445"
$ws.Range("E50").Value = "  +4.31%  "

# Row 51
$ws.Range("D51").Value = "0.9911"
$ws.Range("E51").Value = "  -1.09%  "
